$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.655.87'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.131.25'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.66%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.21'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.58'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.99%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.125.89'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.82%  '
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +14.45%  '
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('E13').Value = '  +4.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.97'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.60%  '
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.649.01'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.18'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.551.36'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.122.45'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '465.03'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.81%  '
$ws.Range('E21').Value = '  +3.40%  '
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.27'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.92%  '
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.96'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +8.40%  '
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.87'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.13'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0879'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +10.80%  '
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.05'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.80%  '
$ws.Range('B36').Value = 'Stacks'
$ws.Range('C36').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.36'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.40'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +13.87%  '
$ws.Range('E38').Value = '  +1.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '51.08'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '452.66'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.76'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.905.44'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('E44').Value = '  +2.48%  '
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('E46').Value = '  +2.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.41'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.49'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  +2.52%  '
